$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Tue Sep 26 21:28:53 EDT 2023"
$ws.Range("B3").Value = "Tue Sep 26 21:29:06 EDT 2023"
$ws.Range("B4").Value = "Tue Sep 26 21:29:19 EDT 2023"
$ws.Range("B5").Value = "Tue Sep 26 21:29:31 EDT 2023"
